# Add a "PETUNJUK" (instructions) column to the student template sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell with a bold style (re-uses/creates the bold font + cellXf automatically).
$ws.Range("G1").Value = "PETUNJUK"
$ws.Range("G1").Font.Bold = $true

# I1 carries the same bold style but stays empty, matching the source template.
$ws.Range("I1").Font.Bold = $true

# Explanatory note placed right under the header.
$ws.Range("G2").Value = "PERHATIAN: Kelas harus sama persis namanya dengan data yang terdaftar sebagai kelas (perhatikan juga huruf besar kecilnya)."

# Widen the new columns so the instructions are readable.
# (Inputs are chosen so the engine's internal rounding of ColumnWidth
# lands on the closest representable width to the authored template.)
$ws.Columns.Item(7).ColumnWidth = 110.33333333333334
$ws.Columns.Item(9).ColumnWidth = 62.66666666666667

# Restore the selection to where the author last left the cursor.
$ws.Range("G8").Select()
